$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.358.21"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "2.963.92"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +10.18%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.52%  "
$ws.Range("D9").Value = "2.959.51"
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.79"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.31"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "3.458.49"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("D18").Value = "2.973.10"
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("D19").Value = "58.463.15"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "420.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.43"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.706"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.73"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  +8.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.71"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +7.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.56"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.02"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0982"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.983"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.99%  "
$ws.Range("D35").Value = "0.0₃0746"
$ws.Range("E35").Value = "  +19.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.70"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.05"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.50"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.66"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.71"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +11.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "395.17"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.81%  "
$ws.Range("D42").Value = "2.729.27"
$ws.Range("E42").Value = "  +4.29%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0345"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.91"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("E47").Value = "  +5.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.99"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.108"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.06"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.66"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +18.23%  "
